$d = $word.ActiveDocument

# Locate the insertion point: right before the "_GoBack" bookmark that
# Word leaves at the very end of the document (inside the last, otherwise
# empty, paragraph).
$bm = $d.Bookmarks.Item("_GoBack")
$insertionPoint = $d.Range($bm.Range.Start, $bm.Range.Start)

# Build the new paragraph's runs as raw WordprocessingML and insert them.
# Wrapping the runs in a <w:p> lets Word merge the run content into the
# existing (target) paragraph at the collapsed insertion point instead of
# creating a brand-new paragraph element.
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$xml = "<w:p $wNs>"
$xml += "<w:r><w:t xml:space=`"preserve`">3. </w:t></w:r>"
$xml += "<w:r><w:t>L</w:t></w:r>"
$xml += "<w:r><w:t xml:space=`"preserve`">e problème qui se trouve lors de l’initialisation de </w:t></w:r>"
$xml += "<w:r><w:t>nombre de livre mensuel autorisé ou à l’ajout du bonus</w:t></w:r>"
$xml += "<w:r><w:t xml:space=`"preserve`"> est que si on veut ajouter d’autre type de package donc on doit toujours changer dans le code de ces deux fonctionnalités. Donc le principe O </w:t></w:r>"
$xml += "<w:r><w:t>(SOLID)</w:t></w:r>"
$xml += "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>"
$xml += "<w:r><w:t>n’est pas respecté.</w:t></w:r>"
$xml += "</w:p>"

$insertionPoint.InsertXML($xml)
